$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 327, pushing existing rows 327-335 down to 328-336.
$ws.Rows.Item(327).Insert()

# Populate the newly inserted row 327 with the new weekly record.
$ws.Cells.Item(327, 1).Value = 7
$ws.Cells.Item(327, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(327, 3).Value = "Ñuble"
$ws.Cells.Item(327, 4).Value = 45075
$ws.Cells.Item(327, 5).Value = 16
$ws.Cells.Item(327, 6).Value = 100112043
$ws.Cells.Item(327, 7).Value = "Pepino ensalada"
$ws.Cells.Item(327, 8).Value = "Sin especificar"
$ws.Cells.Item(327, 9).Value = "Primera"
$ws.Cells.Item(327, 10).Value = 80
$ws.Cells.Item(327, 11).Value = 12000
$ws.Cells.Item(327, 12).Value = 12000
$ws.Cells.Item(327, 13).Value = 12000
$ws.Cells.Item(327, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(327, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(327, 16).Value = 200
$ws.Cells.Item(327, 17).Value = 60
$ws.Cells.Item(327, 18).Value = "Hortaliza"
